# Refresh cryptocurrency Price (column D) and Volume(1h) (column E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.578.71"
$ws.Range("E2").Value = "  +1.60%  "
$ws.Range("D3").Value = "3.440.10"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.34"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.98"
$ws.Range("E6").Value = "  +3.03%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.440.87"
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("E9").Value = "  +8.98%  "
$ws.Range("E10").Value = "  -3.14%  "
$ws.Range("D11").Value = "0.124"
$ws.Range("E11").Value = "  +2.13%  "
$ws.Range("D12").Value = "0.438"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").Value = "4.034.52"
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("E14").Value = "  -2.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000193"
$ws.Range("E15").Value = "  +3.79%  "
$ws.Range("D16").Value = "28.23"
$ws.Range("E16").Value = "  +3.31%  "
$ws.Range("D17").Value = "64.631.61"
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("D18").Value = "3.436.95"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.35"
$ws.Range("E19").Value = "  -1.23%  "
$ws.Range("D20").Value = "14.23"
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("D21").Value = "384.81"
$ws.Range("E21").Value = "  -1.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.16"
$ws.Range("E22").Value = "  -3.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.19"
$ws.Range("E23").Value = "  +1.79%  "
$ws.Range("E24").Value = "  +1.10%  "
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("E26").Value = "  +13.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.73"
$ws.Range("E27").Value = "  +2.55%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.18"
$ws.Range("E30").Value = "  +6.16%  "
$ws.Range("E31").Value = "  +3.43%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.56"
$ws.Range("E33").Value = "  -1.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.60"
$ws.Range("E34").Value = "  -0.35%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  +3.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.83"
$ws.Range("E37").Value = "  +2.91%  "
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("D39").Value = "3.011.12"
$ws.Range("E39").Value = "  +5.14%  "
$ws.Range("E40").Value = "  +1.45%  "
$ws.Range("E41").Value = "  -2.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.30"
$ws.Range("E42").Value = "  -2.71%  "
$ws.Range("D43").Value = "4.54"
$ws.Range("E43").Value = "  +3.68%  "
$ws.Range("E44").Value = "  +2.08%  "
$ws.Range("E45").Value = "  -0.68%  "
$ws.Range("E46").Value = "  +0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.87"
$ws.Range("E47").Value = "  +10.21%  "
$ws.Range("E48").Value = "  -0.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.880"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.18"
$ws.Range("E50").Value = "  +4.61%  "
$ws.Range("E51").Value = "  +3.87%  "
